$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.225.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.427.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '413.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("E9").Value = '  -0.97%  '
$ws.Range("E10").Value = '  +1.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.69'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000217'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.970.95'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.449.34'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.30%  '
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '62.312.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '468.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '90.99'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  +2.83%  '
$ws.Range("E24").Value = '  +4.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.52'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +20.82%  '
$ws.Range("E26").Value = '  +2.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.65'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.43%  '
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("E31").Value = '  -3.89%  '
$ws.Range("E32").Value = '  -2.61%  '
$ws.Range("E33").Value = '  -2.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.40%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0488'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.87%  '
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.05'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.325'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.02%  '
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("E42").Value = '  -1.29%  '
$ws.Range("E43").Value = '  +11.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '145.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.77%  '
$ws.Range("E45").Value = '  +5.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.32'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +18.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.70%  '
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0531'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +31.60%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.20'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.01%  '
